$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = -0.0757
$ws.Range("E2").Value = 0.199
$ws.Range("F2").Value = 0.107
$ws.Range("G2").Value = -0.1073825503355705
$ws.Range("H2").Value = -0.1073825503355705
$ws.Range("I2").Value = -0.1252292748107888
$ws.Range("J2").Value = -0.1252292748107888
$ws.Range("K2").Value = 3020
$ws.Range("L2").Value = 0.6333892617449665
$ws.Range("M2").Value = 1565.2039
$ws.Range("N2").Value = 0.0178454602046773
$ws.Range("O2").Value = 0.5182794370860927
$ws.Range("P2").Value = 139.2039
$ws.Range("Q2").Value = 0.001587114405852092
$ws.Range("R2").Value = 0.04609400662251656
$ws.Range("S2").Value = 1426
$ws.Range("T2").Value = 0.9110634084159898
$ws.Range("U2").Value = 3931
$ws.Range("V2").Value = 0.04481876390966471
$ws.Range("W2").Value = 0.1396983994819132
$ws.Range("X2").Value = 0.1142827812487067
$ws.Range("Y2").Value = 0.0254156182332065
$ws.Range("Z2").Value = 0.2345825685257641
$ws.Range("AA2").Value = -0.02937660493973361
$ws.Range("AB2").Value = 0.1090749198676208
$ws.Range("AC2").Value = -0.1384515248073544
$ws.Range("AD2").Value = 5989
$ws.Range("AE2").Value = 10.46591148920461
$ws.Range("AF2").Value = 5999.465911489205
$ws.Range("AG2").Value = 2068.465911489205
$ws.Range("AH2").Value = 0.06402280367834268
$ws.Range("AI2").Value = 0.1411158715713255
$ws.Range("AJ2").Value = 0.02303997443549341
$ws.Range("AK2").Value = 0.05361016338538074
$ws.Range("AL2").Value = 233
$ws.Range("AM2").Value = 51
$ws.Range("AN2").Value = -13.33853006681514
$ws.Range("AO2").Value = -2.55793991416309
$ws.Range("AP2").Value = -4.60682831066638
$ws.Range("AQ2").Value = -11.68627450980392
$ws.Range("D3").Value = -0.0757
$ws.Range("E3").Value = 0.199
$ws.Range("F3").Value = 0.107
$ws.Range("G3").Value = -0.1073825503355705
$ws.Range("H3").Value = -0.1073825503355705
$ws.Range("I3").Value = -0.1252292748107888
$ws.Range("J3").Value = -0.1252292748107888
$ws.Range("K3").Value = 3020
$ws.Range("L3").Value = 0.6333892617449665
$ws.Range("M3").Value = 1565.2039
$ws.Range("N3").Value = 0.0178454602046773
$ws.Range("O3").Value = 0.5182794370860927
$ws.Range("P3").Value = 139.2039
$ws.Range("Q3").Value = 0.001587114405852092
$ws.Range("R3").Value = 0.04609400662251656
$ws.Range("S3").Value = 1426
$ws.Range("T3").Value = 0.9110634084159898
$ws.Range("U3").Value = 3931
$ws.Range("V3").Value = 0.04481876390966471
$ws.Range("W3").Value = 0.1396983994819132
$ws.Range("X3").Value = 0.1142827812487067
$ws.Range("Y3").Value = 0.0254156182332065
$ws.Range("Z3").Value = 0.2345825685257641
$ws.Range("AA3").Value = -0.02937660493973361
$ws.Range("AB3").Value = 0.1090749198676208
$ws.Range("AC3").Value = -0.1384515248073544
$ws.Range("AD3").Value = 5989
$ws.Range("AE3").Value = 10.46591148920461
$ws.Range("AF3").Value = 5999.465911489205
$ws.Range("AG3").Value = 2068.465911489205
$ws.Range("AH3").Value = 0.06402280367834268
$ws.Range("AI3").Value = 0.1411158715713255
$ws.Range("AJ3").Value = 0.02303997443549341
$ws.Range("AK3").Value = 0.05361016338538074
$ws.Range("AL3").Value = 233
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = -13.33853006681514
$ws.Range("AO3").Value = -2.55793991416309
$ws.Range("AP3").Value = -4.60682831066638
$ws.Range("AQ3").Value = -11.68627450980392